# Apply the "pu" / "pu2" intermediate debug-sheet edit to the FEA workbook.
#
# Summary of the change (per the target diff):
#   - Two new worksheets, "pu" and "pu2", are inserted between "case1" and
#     "origin". Both are duplicates of "case1"'s layout/formatting, but with
#     the force-displacement (E13/E15) inputs replaced by round debug
#     numbers, and rows 30/31 (restraint forces) zeroed out except for a
#     single non-zero entry used to "put a big number" through the solver.
#   - "case1" keeps its data untouched; only its selection/view changes.
#   - "pu2" ends up the active (visible) tab when the file is saved.
#   - case1/pu pick up an explicit paper size + orientation in PageSetup.

$wb = $excel.ActiveWorkbook
$case1 = $wb.Worksheets.Item("case1")

# ------------------------------------------------------------------
# 1. Build "pu" as a duplicate of "case1", positioned right after it.
# ------------------------------------------------------------------
$case1.Copy([System.Reflection.Missing]::Value, $case1)
$pu = $wb.Worksheets.Item(2)
$pu.Name = "pu"

# ------------------------------------------------------------------
# 2. Build "pu2" as a duplicate of "case1", positioned right after "pu".
# ------------------------------------------------------------------
$case1.Copy([System.Reflection.Missing]::Value, $pu)
$pu2 = $wb.Worksheets.Item(3)
$pu2.Name = "pu2"

# ------------------------------------------------------------------
# 3. Overwrite the debug inputs on "pu".
#    E13 (E, steel modulus) and E15 (I, inertia) become flat round
#    numbers instead of the derived case1 formulas; the old restraint
#    force in row 30/31 col G moves to a single 5000 entry in H31.
# ------------------------------------------------------------------
$pu.Range("E13").Formula = "=210000"
$pu.Range("E15").Formula = "=33000000"
$pu.Range("G30").Value = 0
$pu.Range("G31").Value = 0
$pu.Range("H31").Value = 5000

# ------------------------------------------------------------------
# 4. Overwrite the debug inputs on "pu2" - same as "pu" but also puts
#    a big number (400) into H30.
# ------------------------------------------------------------------
$pu2.Range("E13").Formula = "=210000"
$pu2.Range("E15").Formula = "=33000000"
$pu2.Range("G30").Value = 0
$pu2.Range("H30").Value = 400
$pu2.Range("G31").Value = 0
$pu2.Range("H31").Value = 5000

# ------------------------------------------------------------------
# 5. Page setup: case1 and pu get an explicit paper size/orientation.
# ------------------------------------------------------------------
$case1.PageSetup.PaperSize = 9
$case1.PageSetup.Orientation = 1
$pu.PageSetup.PaperSize = 9
$pu.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 6. Selections - restore/update each sheet's view state. "pu2" ends
#    up activated last, so it becomes the workbook's active tab.
# ------------------------------------------------------------------
$case1.Activate()
$case1.Cells.Select()

$pu.Activate()
$pu.Range("H9").Select()
$pu.Cells.Select()

$pu2.Activate()
$pu2.Range("N13").Select()

Write-Host "Sheets now:" ($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "
